$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1829652996845426
$ws.Range("C2").Value = 0.5772870662460567
$ws.Range("J2").Value = 0.009463722397476341
$ws.Range("P2").Value = 0.1545741324921136
$ws.Range("S2").Value = 0.07570977917981073
$ws.Range("B3").Value = 0.005154639175257732
$ws.Range("C3").Value = 0.05154639175257732
$ws.Range("J3").Value = 0.03092783505154639
$ws.Range("P3").Value = 0.7268041237113402
$ws.Range("S3").Value = 0.1855670103092784
$ws.Range("J4").Value = 0.02222222222222222
$ws.Range("P4").Value = 0.6444444444444445
$ws.Range("S4").Value = 0.3333333333333333
$ws.Range("B6").Value = 0.05825242718446602
$ws.Range("D6").Value = 0.004854368932038835
$ws.Range("F6").Value = 0.04854368932038835
$ws.Range("J6").Value = 0.2718446601941747
$ws.Range("O6").Value = 0.02427184466019417
$ws.Range("Q6").Value = 0.1893203883495146
$ws.Range("R6").Value = 0.04854368932038835
$ws.Range("S6").Value = 0.354368932038835
$ws.Range("B7").Value = 0.1125541125541126
$ws.Range("D7").Value = 0.0303030303030303
$ws.Range("F7").Value = 0.06060606060606061
$ws.Range("J7").Value = 0.1255411255411255
$ws.Range("O7").Value = 0.01731601731601732
$ws.Range("Q7").Value = 0.1471861471861472
$ws.Range("R7").Value = 0.08658008658008658
$ws.Range("S7").Value = 0.4199134199134199
$ws.Range("B8").Value = 0.1002178649237473
$ws.Range("D8").Value = 0.01525054466230937
$ws.Range("F8").Value = 0.05228758169934641
$ws.Range("J8").Value = 0.1241830065359477
$ws.Range("O8").Value = 0.0130718954248366
$ws.Range("Q8").Value = 0.2069716775599129
$ws.Range("R8").Value = 0.06100217864923747
$ws.Range("S8").Value = 0.4270152505446623
$ws.Range("B9").Value = 0.101063829787234
$ws.Range("F9").Value = 0.05851063829787234
$ws.Range("J9").Value = 0.1117021276595745
$ws.Range("Q9").Value = 0.2287234042553191
$ws.Range("R9").Value = 0.1063829787234043
$ws.Range("S9").Value = 0.3936170212765958
$ws.Range("B10").Value = 0.1115214180206795
$ws.Range("D10").Value = 0.02141802067946824
$ws.Range("F10").Value = 0.06277695716395865
$ws.Range("J10").Value = 0.1344165435745938
$ws.Range("O10").Value = 0.01698670605612999
$ws.Range("Q10").Value = 0.2518463810930576
$ws.Range("R10").Value = 0.07090103397341212
$ws.Range("S10").Value = 0.3301329394387001
$ws.Range("G11").Value = 0.1529745042492918
$ws.Range("J11").Value = 0.06515580736543909
$ws.Range("K11").Value = 0.2011331444759207
$ws.Range("L11").Value = 0.5637393767705382
$ws.Range("S11").Value = 0.0169971671388102
$ws.Range("G12").Value = 0.7826086956521739
$ws.Range("J12").Value = 0.1497584541062802
$ws.Range("K12").Value = 0.004830917874396135
$ws.Range("L12").Value = 0.04347826086956522
$ws.Range("S12").Value = 0.01932367149758454
$ws.Range("F13").Value = 0.02631578947368421
$ws.Range("G13").Value = 0.5789473684210527
$ws.Range("J13").Value = 0.3421052631578947
$ws.Range("S13").Value = 0.05263157894736842
$ws.Range("F15").Value = 0.00881057268722467
$ws.Range("H15").Value = 0.1497797356828194
$ws.Range("I15").Value = 0.06607929515418502
$ws.Range("J15").Value = 0.3832599118942731
$ws.Range("K15").Value = 0.06607929515418502
$ws.Range("M15").Value = 0.01762114537444934
$ws.Range("O15").Value = 0.05286343612334802
$ws.Range("S15").Value = 0.2555066079295154
$ws.Range("F16").Value = 0.01395348837209302
$ws.Range("H16").Value = 0.186046511627907
$ws.Range("I16").Value = 0.09302325581395349
$ws.Range("J16").Value = 0.386046511627907
$ws.Range("K16").Value = 0.1255813953488372
$ws.Range("M16").Value = 0.009302325581395349
$ws.Range("N16").Value = 0.004651162790697674
$ws.Range("O16").Value = 0.05116279069767442
$ws.Range("S16").Value = 0.1302325581395349
$ws.Range("F17").Value = 0.007312614259597806
$ws.Range("H17").Value = 0.1425959780621572
$ws.Range("I17").Value = 0.08409506398537477
$ws.Range("J17").Value = 0.4552102376599634
$ws.Range("K17").Value = 0.1005484460694698
$ws.Range("M17").Value = 0.02010968921389397
$ws.Range("O17").Value = 0.05484460694698354
$ws.Range("S17").Value = 0.1352833638025594
$ws.Range("F18").Value = 0.01142857142857143
$ws.Range("H18").Value = 0.1714285714285714
$ws.Range("I18").Value = 0.1085714285714286
$ws.Range("J18").Value = 0.44
$ws.Range("K18").Value = 0.07428571428571429
$ws.Range("M18").Value = 0.01142857142857143
$ws.Range("O18").Value = 0.05142857142857143
$ws.Range("S18").Value = 0.1314285714285714
$ws.Range("F19").Value = 0.01931993817619784
$ws.Range("H19").Value = 0.2202472952086553
$ws.Range("I19").Value = 0.07032457496136012
$ws.Range("J19").Value = 0.3469860896445132
$ws.Range("K19").Value = 0.1313755795981453
$ws.Range("M19").Value = 0.0170015455950541
$ws.Range("N19").Value = 0.0007727975270479134
$ws.Range("O19").Value = 0.07573415765069552
$ws.Range("S19").Value = 0.1182380216383308
